$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://tamian.uk/contact-us/"

for ($i = 1; $i -le 120; $i++) {
    $ws.Range("N$i").Value = $url
}
